# add: creation big label and input protection
#
# 1) Rename existing sheet "Sheet1" -> "etiqueta_peq"
# 2) Add a new sheet "etiqueta_grande" right after it
# 3) Fix column A (rows 2-13) on "etiqueta_peq": "s" -> "32" (kept as text)
# 4) Populate "etiqueta_grande" with its header row (copied formatting from
#    the "etiqueta_peq" header) and its 4 data rows

$wb = $excel.ActiveWorkbook

# --- 1) rename the first/original sheet ---------------------------------
$wsPeq = $wb.Worksheets.Item(1)
$wsPeq.Name = "etiqueta_peq"

# --- 2) add the new sheet right after it ---------------------------------
$wsGrande = $wb.Worksheets.Add($null, $wsPeq)
$wsGrande.Name = "etiqueta_grande"

# --- 3) fix "etiqueta_peq" column A values (rows 2-13): "s" -> "32" -----
# Written with a leading apostrophe so it lands as text (matches the
# original inline-string "32", not a number), then the style is reset
# back to Normal so no stray number-format/quote-prefix style sticks.
for ($r = 2; $r -le 13; $r++) {
    $cell = $wsPeq.Cells.Item($r, 1)
    $cell.Value = "'32"
    $cell.Style = "Normal"
}

# --- 4) populate "etiqueta_grande" ---------------------------------------

# Header formatting: reuse the exact header style from "etiqueta_peq"
# (bold, bordered, centered) by copy/paste-special of formats only.
$wsPeq.Range("A1:F1").Copy() | Out-Null
$wsGrande.Range("A1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$headers = @("CLIENTE", "DESTINO", "PROVEEDOR", "OC", "NRO. DE GUIA", "ASN", "CANT BULTOS", "PESO", "LPN", "TIPO")
for ($c = 1; $c -le $headers.Length; $c++) {
    $wsGrande.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Data rows 2-5. Columns D, E, F, G are text (even though numeric-looking),
# column H is a genuine number; everything else is plain text.
$rows = @(
    @("COMPAÑIA MINERA DOÑA INES DE COLLAHUASI", "BODEGA ROSARIO", "COMERCIAL, SERVICIOS E INGENIERIA CSI SPA", "32", "3", "2", "01 DE 03", 2, "SAL0000004491", "BULTO"),
    @("COMPAÑIA MINERA DOÑA INES DE COLLAHUASI", "BODEGA ROSARIO", "COMERCIAL, SERVICIOS E INGENIERIA CSI SPA", "32", "3", "2", "02 DE 03", 22, "SAL0000004493", "BULTO"),
    @("COMPAÑIA MINERA DOÑA INES DE COLLAHUASI", "BODEGA ROSARIO", "COMERCIAL, SERVICIOS E INGENIERIA CSI SPA", "32", "3", "2", "03 DE 03", 3, "SAL0000004528", "BULTO"),
    @("COMPAÑIA MINERA DOÑA INES DE COLLAHUASI", "BODEGA ROSARIO", "COMERCIAL, SERVICIOS E INGENIERIA CSI SPA", "32", "3", "2", "01 DE 01", 2, "Pallet1", "PALLET")
)

$textCols = @(1, 2, 3, 4, 5, 6, 7, 9, 10)
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $cell = $wsGrande.Cells.Item($r, $c)
        if ($textCols -contains $c) {
            $cell.Value = "'" + $rowData[$c - 1]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $rowData[$c - 1]
        }
    }
}

$wsPeq.Select()
